$d = $word.ActiveDocument

function Find-ParagraphContaining($doc, [string]$needle) {
    foreach ($para in $doc.Paragraphs) {
        if ($para.Range.Text -like "*$needle*") {
            return $para
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Change 1: "Kommandos:" -> "Kommandos für Website:" (in the heading paragraph
# right after the title), with a bookmark ("_GoBack") placed right before
# "Website" (it used to sit inside "Impressum" further down the document).
# ---------------------------------------------------------------------------

$kommandosPara = Find-ParagraphContaining $d "Kommandos:"
if ($null -eq $kommandosPara) {
    throw "Could not find the 'Kommandos:' paragraph."
}

$kommandosRange = $kommandosPara.Range.Duplicate
$kommandosRange.Find.ClearFormatting()
$found = $kommandosRange.Find.Execute("Kommandos:", $false, $false, $false, $false, $false, `
                                       $true, 1, $false, "Kommandos für Website:", 2)
if (-not $found) {
    throw "Could not replace 'Kommandos:' with 'Kommandos für Website:'."
}

# Re-locate the (now replaced) paragraph so we can find exactly where
# "Website" starts and drop the bookmark immediately in front of it.
$kommandosPara2 = Find-ParagraphContaining $d "Kommandos für Website:"
$websiteRange = $kommandosPara2.Range.Duplicate
$websiteRange.Find.ClearFormatting()
$websiteRange.Find.Execute("Website:") | Out-Null
$bookmarkSpot = $d.Range($websiteRange.Start, $websiteRange.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)

# ---------------------------------------------------------------------------
# Change 2: "Impre" + bookmark + "ssum" -> plain "Impressum" run (the old
# "_GoBack" bookmark that used to live inside this word is gone now, it
# moved to the spot above).
# ---------------------------------------------------------------------------

$impressumPara = Find-ParagraphContaining $d "Impressum"
if ($null -eq $impressumPara) {
    throw "Could not find the 'Impressum' paragraph."
}

$impressumRange = $impressumPara.Range.Duplicate
$impressumRange.Find.ClearFormatting()
$impressumRange.Find.Execute("Impressum", $false, $false, $false, $false, $false, `
                              $true, 1, $false, "Impressum", 2) | Out-Null
